# Fruta / hortaliza, semanal
# Weekly refresh: insert the newest week's two rows at the top of the data
# block (pushing the existing rows down by two) and update a handful of
# values on the rows that already existed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Make room for the two new rows of data (new rows 2 and 3). This
#    shifts the previous rows 2-9 down to rows 4-11.
# ---------------------------------------------------------------------
$ws.Rows("2:3").Insert()

# The blank rows inherit the header's bold style from Insert(); strip
# that back to the default (unstyled) look used by the other data rows.
$ws.Rows("2:3").ClearFormats()

# Column D carries a date number format in every data row; re-apply it
# to the two new date cells.
$ws.Range("D2:D3").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# ---------------------------------------------------------------------
# 2) Populate the two new rows (common columns + the 2023-04-27 values).
# ---------------------------------------------------------------------
function Set-DataRow($row, $date, $quality, $volume, $minPrice, $maxPrice, $avgPrice, $unit, $origin, $pricePerKg, $kgPerUnit) {
    $ws.Cells.Item($row, 1).Value2 = 7
    $ws.Cells.Item($row, 2).Value2 = "Terminal Hortofrutícola Agro Chillán"
    $ws.Cells.Item($row, 3).Value2 = "Ñuble"
    $ws.Cells.Item($row, 4).Value2 = $date
    $ws.Cells.Item($row, 5).Value2 = 16
    $ws.Cells.Item($row, 6).Value2 = "Fruta"
    $ws.Cells.Item($row, 7).Value2 = 100104
    $ws.Cells.Item($row, 8).Value2 = "Frutos de pepita"
    $ws.Cells.Item($row, 9).Value2 = 100104003
    $ws.Cells.Item($row, 10).Value2 = "Membrillo"
    $ws.Cells.Item($row, 11).Value2 = "Champion"
    $ws.Cells.Item($row, 12).Value2 = $quality
    $ws.Cells.Item($row, 13).Value2 = $volume
    $ws.Cells.Item($row, 14).Value2 = $minPrice
    $ws.Cells.Item($row, 15).Value2 = $maxPrice
    $ws.Cells.Item($row, 16).Value2 = $avgPrice
    $ws.Cells.Item($row, 17).Value2 = $unit
    $ws.Cells.Item($row, 18).Value2 = $origin
    $ws.Cells.Item($row, 19).Value2 = $pricePerKg
    $ws.Cells.Item($row, 20).Value2 = $kgPerUnit
}

Set-DataRow 2 45043 "Especial" 40 13000 13000 13000 "`$/caja 18 kilos empedrada" "Región de O'Higgins" 722 18
Set-DataRow 3 45043 "Primera"  50 12000 12000 12000 "`$/caja 18 kilos empedrada" "Región de O'Higgins" 667 18

# ---------------------------------------------------------------------
# 3) Update the rows that were pushed down (now rows 4-11) with their
#    new values.
# ---------------------------------------------------------------------
Set-DataRow 4  45033 "Especial" 60  13000 13000 13000 "`$/caja 18 kilos empedrada" "Región de O'Higgins" 722 18
Set-DataRow 5  45033 "Primera"  80  12000 12000 12000 "`$/caja 18 kilos empedrada" "Región de O'Higgins" 667 18
Set-DataRow 6  45020 "Primera"  60  12000 12000 12000 "`$/caja 18 kilos granel"    "Región de O'Higgins" 667 18
Set-DataRow 7  45021 "Primera"  50  12000 12000 12000 "`$/caja 18 kilos granel"    "Región de O'Higgins" 667 18
Set-DataRow 8  44699 "Especial" 60  13000 13000 13000 "`$/caja 15 kilos granel"    "Provincia de Curicó" 867 15
Set-DataRow 9  44699 "Primera"  120 11000 12000 11500 "`$/caja 15 kilos granel"    "Provincia de Curicó" 767 15
Set-DataRow 10 45040 "Especial" 50  13000 13000 13000 "`$/caja 18 kilos empedrada" "Región de O'Higgins" 722 18
Set-DataRow 11 45040 "Primera"  40  12000 12000 12000 "`$/caja 18 kilos empedrada" "Región de O'Higgins" 667 18

Write-Host "Applied weekly Membrillo update"
